$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

$ws.Range("F2").Value = "2021-10-05 13:41:48.246085"
$ws.Range("F3").Value = "2021-10-05 13:41:48.246096"
$ws.Range("F4").Value = "2021-10-05 13:41:48.246100"
$ws.Range("F5").Value = "2021-10-05 13:41:48.246103"
$ws.Range("F6").Value = "2021-10-05 13:41:48.246107"
